# Add two new columns "I0" (col I) and "IF" (col J) to the sheet,
# matching the header style used by the existing columns (e.g. H1),
# and fill in the numeric values for rows 2-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy the style of the existing last header cell (H1) to I1/J1
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-21 for column I ("I0") and column J ("IF")
$dataI = @(8,8,7,5,5,7,8,6,8,7,4,8,9,7,8,6,6,6,7,6)
$dataJ = @(9,8,7,5,6,7,9,6,8,7,5,8,9,7,8,6,6,6,7,6)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
